$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - reuse the same formatting as the other header cells
# (e.g. G1: bold, thin border, centered horizontally, top vertical align)
# by copying its format onto H1, then set the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New "Save" data column for rows 2-5 (values taken from the diff)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
